$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tmp = $ws.Range("Z1")

$ws.Range("D2").Value = "40.086.18"
$ws.Range("E2").Value = "  +2.80%  "
$ws.Range("D3").Value = "2.237.86"
$ws.Range("E3").Value = "  +1.28%  "
$ws.Range("E4").Value = "  +0.08%  "
$tmp.Value = "'294.73"
$tmp.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = "  -0.34%  "
$tmp.Value = "'86.31"
$tmp.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = "  +7.70%  "
$ws.Range("E7").Value = "  +2.35%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +3.65%  "
$tmp.Value = "'31.16"
$tmp.Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = "  +11.63%  "
$tmp.Value = "'0.0792"
$tmp.Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = "  +2.69%  "
$tmp.Value = "'47.02"
$tmp.Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Value = "  +2.79%  "
$ws.Range("E13").Value = "  +1.09%  "
$tmp.Value = "'6.46"
$tmp.Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = "  +6.32%  "
$ws.Range("D15").Value = "2.580.75"
$ws.Range("E15").Value = "  +0.95%  "
$tmp.Value = "'14.17"
$tmp.Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("E16").Value = "  +2.29%  "
$ws.Range("D17").Value = "2.187.06"
$ws.Range("E17").Value = "  -2.09%  "
$ws.Range("D19").Value = "40.019.67"
$ws.Range("E19").Value = "  +2.94%  "
$ws.Range("D20").Value = "0.0₃0891"
$ws.Range("E20").Value = "  +4.03%  "
$ws.Range("E21").Value = "  +2.33%  "
$tmp.Value = "'10.82"
$tmp.Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value = "  +10.83%  "
$tmp.Value = "'65.41"
$tmp.Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = "  +1.31%  "
$tmp.Value = "'234.89"
$tmp.Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = "  +4.71%  "
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("E26").Value = "  +3.69%  "
$ws.Range("E27").Value = "  +6.76%  "
$ws.Range("E28").Value = "  +2.91%  "
$ws.Range("E29").Value = "  +3.38%  "
$tmp.Value = "'9.23"
$tmp.Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = "  +4.21%  "
$tmp.Value = "'33.37"
$tmp.Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = "  +7.65%  "
$tmp.Value = "'152.35"
$tmp.Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = "  +2.88%  "
$ws.Range("E33").Value = "  +0.03%  "
$tmp.Value = "'4.88"
$tmp.Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value = "  +3.48%  "
$tmp.Value = "'0.0722"
$tmp.Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = "  +5.36%  "
$tmp.Value = "'2.39"
$tmp.Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = "  +2.53%  "
$tmp.Value = "'16.35"
$tmp.Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = "  +14.83%  "
$ws.Range("E38").Value = "  +3.34%  "
$ws.Range("E39").Value = "  +5.36%  "
$ws.Range("E40").Value = "  +4.26%  "
$ws.Range("E41").Value = "  +7.38%  "
$ws.Range("E42").Value = "  +7.04%  "
$ws.Range("D43").Value = "2.052.76"
$ws.Range("E43").Value = "  +8.51%  "
$tmp.Value = "'2.24"
$tmp.Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = "  +7.69%  "
$ws.Range("E45").Value = "  +6.75%  "
$tmp.Value = "'10.02"
$tmp.Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = "  +14.00%  "
$tmp.Value = "'16.30"
$tmp.Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = "  +1.84%  "
$tmp.Value = "'2.57"
$tmp.Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = "  +2.88%  "
$ws.Range("D49").Value = "2.451.52"
$ws.Range("E49").Value = "  +1.06%  "
$tmp.Value = "'70.94"
$tmp.Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value = "  +1.66%  "
$ws.Range("E51").Value = "  +14.96%  "
$tmp.Clear()
